$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price column cells that are being updated,
# to avoid Excel auto-converting numeric-looking strings to numbers.
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D11", "D13", "D14", "D15", "D16", "D17", "D20", "D23", "D24", "D25", "D27", "D28", "D29", "D31", "D32", "D33", "D34", "D35", "D38", "D39", "D40", "D41", "D43", "D44", "D45", "D46", "D48", "D50", "D51")
foreach ($c in $priceCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.768.52"
$ws.Range("E2").Value = "  +0.08%  "

$ws.Range("D3").Value = "3.797.78"
$ws.Range("E3").Value = "  +0.52%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "600.87"
$ws.Range("E5").Value = "  +0.92%  "

$ws.Range("D6").Value = "165.21"
$ws.Range("E6").Value = "  -1.02%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  -0.51%  "

$ws.Range("E9").Value = "  -0.09%  "

$ws.Range("E10").Value = "  +0.75%  "

$ws.Range("D11").Value = "6.39"
$ws.Range("E11").Value = "  +1.63%  "

$ws.Range("E12").Value = "  -1.29%  "

$ws.Range("D13").Value = "35.78"
$ws.Range("E13").Value = "  -0.69%  "

$ws.Range("D14").Value = "4.434.12"
$ws.Range("E14").Value = "  +0.45%  "

$ws.Range("D15").Value = "3.796.28"
$ws.Range("E15").Value = "  +0.49%  "

$ws.Range("D16").Value = "67.787.33"
$ws.Range("E16").Value = "  +0.21%  "

$ws.Range("D17").Value = "18.41"
$ws.Range("E17").Value = "  -0.65%  "

$ws.Range("E18").Value = "  +1.39%  "

$ws.Range("E19").Value = "  +0.49%  "

$ws.Range("D20").Value = "463.72"
$ws.Range("E20").Value = "  +0.91%  "

$ws.Range("E21").Value = "  -2.50%  "

$ws.Range("E22").Value = "  +0.45%  "

$ws.Range("D23").Value = "0.0000147"
$ws.Range("E23").Value = "  -5.23%  "

$ws.Range("D24").Value = "82.82"
$ws.Range("E24").Value = "  -0.69%  "

$ws.Range("D25").Value = "12.05"
$ws.Range("E25").Value = "  +0.53%  "

$ws.Range("E26").Value = "  -0.62%  "

$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").Value = "10.01"
$ws.Range("E27").Value = "  -0.04%  "

$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.04%  "

$ws.Range("D29").Value = "3.947.00"
$ws.Range("E29").Value = "  +0.53%  "

$ws.Range("E30").Value = "  -0.15%  "

$ws.Range("D31").Value = "7.41"
$ws.Range("E31").Value = "  +2.64%  "

$ws.Range("D32").Value = "2.20"
$ws.Range("E32").Value = "  -1.51%  "

$ws.Range("D33").Value = "29.24"
$ws.Range("E33").Value = "  -1.02%  "

$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.01%  "

$ws.Range("D35").Value = "9.03"
$ws.Range("E35").Value = "  -0.49%  "

$ws.Range("E36").Value = "  -0.56%  "

$ws.Range("E37").Value = "  +0.75%  "

$ws.Range("D38").Value = "0.998"
$ws.Range("E38").Value = "  +0.31%  "

$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").Value = "3.23"
$ws.Range("E39").Value = "  -4.23%  "

$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").Value = "5.76"
$ws.Range("E40").Value = "  -0.15%  "

$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  +0.04%  "

$ws.Range("D43").Value = "45.27"
$ws.Range("E43").Value = "  -0.55%  "

$ws.Range("D44").Value = "47.65"
$ws.Range("E44").Value = "  -0.85%  "

$ws.Range("D45").Value = "0.298"
$ws.Range("E45").Value = "  -0.30%  "

$ws.Range("D46").Value = "151.35"
$ws.Range("E46").Value = "  +0.99%  "

$ws.Range("E47").Value = "  +11.18%  "

$ws.Range("D48").Value = "27.51"
$ws.Range("E48").Value = "  +3.23%  "

$ws.Range("E49").Value = "  +0.36%  "

$ws.Range("D50").Value = "395.21"
$ws.Range("E50").Value = "  +0.38%  "

$ws.Range("D51").Value = "1.85"
$ws.Range("E51").Value = "  +1.80%  "

# Reset style index back to default (Normal) for the forced-text cells
# while keeping their text number format / values intact.
foreach ($c in $priceCells) {
    $ws.Range($c).Style = "Normal"
}